# Apply the crypto price/volume refresh described by the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.436.50"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.673.74"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.82%  "
$ws.Range("D5").Value = "'221.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").Value = "'0.5349"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").Value = "'0.2670"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("D9").Value = "'0.06409"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("D11").Value = "'0.07847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.701.62"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.553"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "1.903.72"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "'0.5656"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "0.0₅8198"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'66.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "26.466.68"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").Value = "'4.739"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").Value = "'198.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.86%  "
$ws.Range("D22").Value = "'10.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").Value = "'6.076"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Value = "'146.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "'7.261"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "'16.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("D29").Value = "'1.504"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").Value = "'0.05888"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").Value = "'1.289"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").Value = "'3.587"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "'3.316"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").Value = "'1.618"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("D35").Value = "'0.9714"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("D36").Value = "'2.855"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "1.082.37"
$ws.Range("E40").Value = "  +3.94%  "
$ws.Range("D41").Value = "'5.928"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("D42").Value = "'0.8668"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "'104.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "1.813.58"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").Value = "'58.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("E47").Value = "  -4.16%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "'0.4418"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("D50").Value = "'8.056"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").Value = "'0.05172"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
